$d = $word.ActiveDocument

# --- Step 1: insert the new narrative paragraphs right before the
#     paragraph that carries the "_GoBack" bookmark. Word keeps the
#     bookmark anchored to its original (now pushed-down) paragraph
#     when text + paragraph breaks are inserted immediately before it,
#     which matches the target structure (a run of brand-new <w:p>
#     elements followed by the original, untouched bookmark paragraph).
$bm = $d.Bookmarks.Item("_GoBack")
$ins = $bm.Range.Duplicate
$ins.Collapse(1)

$newText = "Una vez conocidos los datos con los que trataremos, observaremos como se representan estos en su totalidad y que relación guardan con otros atributos, con el fin de eliminar atributos redundantes y establecer relaciones para predecir una variable objetivo la cual tendremos que especificar también a partir de la representación de estos datos.`rPara comenzar a observar cómo se representan los datos, usaremos funciones de la librería pandas primeramente.`rEl primer paso es saber la dimensión total de los datos obtenidos en nuestra variable dataset el cual podremos saber con la función shape(), así obtenemos que nuestro dataset contiene un total de 4177 filas.`r`rPosteriormente necesitamos saber si existen valores nulos en alguna de las filas de nuestro dataset debido a que estos valores no nos servirán posteriormente para el análisis de datos, así que se hace uso de la función isnull combinada con la función sum() de manera que obtenemos los valores nulos que existen por cada columna de nuestro dataset, dando así como resultado que no contiene valores nulos, de manera que en este apartado no descartaremos datos todavía.`r`rDespués de haber realizado estos dos simples pasos, pasamos a utilizar la función describe() la cual nos da datos relevantes sobre cada columna de nuestro dataset, tal como la media de valores, los máximos, los mínimos, los percentiles, etc.`rDe manera que así podemos observar posibles inconsistencias en los datos para seguir realizando la criba de estos. En este caso vemos que los valores representados con estas funciones parecen correctos a priori, pero nos llama la atención que la altura mínima encontrada es igual a 0, cosa que no tiene sentido, así que se procede a la búsqueda de las filas del dataset que contengan este valor filtrando por el atributo “Height” e igualándolo a 0.`rEsta consulta nos devuelve dos filas, las cuales habrá que descartar ya que aportan inconsistencias. Para descartar estas filas hacemos uso de la función drop() a la cual como parámetros le pasaremos en una lista los índices de las filas con altura igual a 0.`rPosteriormente se comprueba de nuevo la misma consulta para ver que no encuentra resultados y se comprueba de nuevo la dimensión del dataset para ver que contiene dos filas menos.`r`rSiguiendo estos procedimientos, sabemos que se puede calcular la edad de un molusco sumando 1,5 al número de sus anillos, así que se crea una variable llamada age que sea igual al resultado de esta suma, utilizando la siguiente sentencia: `rdataset[‘age’] = dataset.Rings +1.5`r`rComo paso extra, podemos decir que renombramos el nombre de las variables para evitar problemas luego de consultas como pueden ocurrir a partir de espacios entre los nombres, o el uso de mayúsculas. Así que se cambian los espacios por “_”, y se sustituyen las mayúsculas por minúsculas usando la función rename de la librería pandas pasando por parámetro las columnas con su nombre antiguo y especificando el nuevo.`rPara comprobar que estos cambios han sido ejecutados correctamente podremos utilizar de nuevo la función head() para ver mismamente los primeros nuevos valores y así comprobar que los datos han sido actualizados y además aparece el nuevo atributo “age”.`r `r"

$ins.InsertBefore($newText)

# --- Step 2: bold the “age” attribute name inside the paragraph that
#     ends with: ...además aparece el nuevo atributo "age".
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("“age”.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $boldRange = $findRange.Duplicate
    $boldRange.SetRange($findRange.Start, $findRange.End - 1)
    $boldRange.Font.Bold = 1
}

# --- Step 3: remove the two now-superfluous empty paragraphs that sit
#     right after the bookmark paragraph (they previously carried bold
#     paragraph-mark formatting and are dropped in this revision).
#     Locate the bookmark paragraph's numeric index robustly (its
#     position shifted because of the paragraphs inserted above).
$bm2 = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm2.Range.Start
$bmIndex = -1
$idx = 1
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -le $bmStart -and $para.Range.End -ge $bmStart) {
        $bmIndex = $idx
    }
    $idx = $idx + 1
}

# The paragraph immediately after the bookmark paragraph keeps its
# (non-bold) formatting. Of the next four, the 1st and 3rd carry bold
# paragraph-mark formatting and are the ones dropped in this revision
# (their non-bold siblings at relative +1/+3 before each deletion
# survive), so delete relative position +2, then (after the shift)
# the new relative position +3.
$pDel1 = $d.Paragraphs.Item($bmIndex + 2)
$pDel1.Range.Delete()
$pDel2 = $d.Paragraphs.Item($bmIndex + 3)
$pDel2.Range.Delete()
